$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells on row 1
$ws.Range("C1").Value = "Phần trăm lãi suất"
$ws.Range("D1").Value = 40

# Update header on row 2
$ws.Range("D2").Value = "Giá bán"

# Append numbering suffixes to existing book titles (columns B3:B8)
$ws.Range("B3").Value = "Conan 1"
$ws.Range("B4").Value = "Năm mươi Sắc thái 1"
$ws.Range("B5").Value = "Cho tôi một vé đi tuổi thơ 2"
$ws.Range("B6").Value = "Chiến Thắng Con Quỷ Trong Bạn 3"
$ws.Range("B7").Value = "Đất Rừng Phương Nam 4"
$ws.Range("B8").Value = "Harry Potter Và Hòn Đá Phù Thuỷ 4"

# Clear out the remaining auto-filled supplier/date rows
$ws.Range("A9:D22").ClearContents()

# Restore the active selection used when the workbook was last saved
$ws.Range("N2").Select()
